$wb = $excel.ActiveWorkbook

# --- Update the "Bingo" password values on the approver lookup sheets ---
# These cells (B2:B6) are hyperlinked mailto cells whose display text moves
# from "Bingo@1234567" to "Bingo@12345".
$ws5 = $wb.Worksheets.Item("FirstLevelApprover")
$ws5.Range("B2").Value = "Bingo@12345"
$ws5.Range("B3").Value = "Bingo@12345"
$ws5.Range("B4").Value = "Bingo@12345"
$ws5.Range("B5").Value = "Bingo@12345"
$ws5.Range("B6").Value = "Bingo@12345"

$ws6 = $wb.Worksheets.Item("Approver")
$ws6.Range("B2").Value = "Bingo@12345"
$ws6.Range("B3").Value = "Bingo@12345"
$ws6.Range("B4").Value = "Bingo@12345"
$ws6.Range("B5").Value = "Bingo@12345"
$ws6.Range("B6").Value = "Bingo@12345"

# --- Update selections / active sheet state to match the saved view ---
$ws5.Activate()
[void]$ws5.Range("B7").Select()

$ws6.Activate()
[void]$ws6.Range("H22").Select()
